$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 (Home) target depth numbers ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 428
$wsOff.Range("C2").Value = 320
$wsOff.Range("D2").Value = 120
$wsOff.Range("E2").Value = 64
$wsOff.Range("F2").Value = 6

# --- DEF sheet: row 2 (Home) target depth numbers ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 459
$wsDef.Range("C2").Value = 320
$wsDef.Range("D2").Value = 117
$wsDef.Range("E2").Value = 45
